$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet
$ws.Name = "ValidLogin"

# Set the new header / credential values
$ws.Range("A2").Value = "Admin"
$ws.Range("B2").Value = "admin123"
$ws.Range("A1").Value = "Username"
$ws.Range("B1").Value = "Password"

# Move / update selection to mirror the target workbook (row 7 selected)
$ws.Range("A7:XFD7").Select()
